$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh adds two newer observations (11-Jul-2022) at the top of
# the data block and pushes the existing rows 118-125 down to 120-127.
$ws.Range("A118:A119").EntireRow.Insert()

# New row 118: Membrillo Champion "Especial"
$ws.Cells.Item(118, 1).Value2 = 6
$ws.Cells.Item(118, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(118, 3).Value2 = "Metropolitana"
$ws.Cells.Item(118, 4).Value2 = 44753
$ws.Cells.Item(118, 5).Value2 = 13
$ws.Cells.Item(118, 6).Value2 = "Fruta"
$ws.Cells.Item(118, 7).Value2 = 100104
$ws.Cells.Item(118, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(118, 9).Value2 = 100104003
$ws.Cells.Item(118, 10).Value2 = "Membrillo"
$ws.Cells.Item(118, 11).Value2 = "Champion"
$ws.Cells.Item(118, 12).Value2 = "Especial"
$ws.Cells.Item(118, 13).Value2 = 8
$ws.Cells.Item(118, 14).Value2 = 280000
$ws.Cells.Item(118, 15).Value2 = 280000
$ws.Cells.Item(118, 16).Value2 = 280000
$ws.Cells.Item(118, 17).Value2 = "`$/bins (450 kilos)"
$ws.Cells.Item(118, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(118, 19).Value2 = 622
$ws.Cells.Item(118, 20).Value2 = 450

# New row 119: Membrillo Champion "Segunda"
$ws.Cells.Item(119, 1).Value2 = 6
$ws.Cells.Item(119, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(119, 3).Value2 = "Metropolitana"
$ws.Cells.Item(119, 4).Value2 = 44753
$ws.Cells.Item(119, 5).Value2 = 13
$ws.Cells.Item(119, 6).Value2 = "Fruta"
$ws.Cells.Item(119, 7).Value2 = 100104
$ws.Cells.Item(119, 8).Value2 = "Frutos de pepita"
$ws.Cells.Item(119, 9).Value2 = 100104003
$ws.Cells.Item(119, 10).Value2 = "Membrillo"
$ws.Cells.Item(119, 11).Value2 = "Champion"
$ws.Cells.Item(119, 12).Value2 = "Segunda"
$ws.Cells.Item(119, 13).Value2 = 10
$ws.Cells.Item(119, 14).Value2 = 200000
$ws.Cells.Item(119, 15).Value2 = 200000
$ws.Cells.Item(119, 16).Value2 = 200000
$ws.Cells.Item(119, 17).Value2 = "`$/bins (450 kilos)"
$ws.Cells.Item(119, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(119, 19).Value2 = 444
$ws.Cells.Item(119, 20).Value2 = 450
